$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement rows (16-19) list two workers (CC 73116093 / ARMANDO
# ENRIQUE JULIAO BLANCO and CC 73114232 / ARMANDO RAFAEL MORELOS ALVIS),
# each with two overdue periods (1801 and 1802). The database update
# re-sorts the rows so they are grouped by period (1801 then 1802) instead
# of by worker, while keeping each worker/period combination intact.

$ws.Range("E16").Value = "1801"

$ws.Range("C17").Value = "73114232"
$ws.Range("D17").Value = "ARMANDO RAFAEL MORELOS ALVIS"
$ws.Range("E17").Value = "1801"

$ws.Range("C18").Value = "73116093"
$ws.Range("D18").Value = "ARMANDO ENRIQUE JULIAO BLANCO"
$ws.Range("E18").Value = "1802"

$ws.Range("C19").Value = "73114232"
$ws.Range("D19").Value = "ARMANDO RAFAEL MORELOS ALVIS"
$ws.Range("E19").Value = "1802"
